# Commit: "Tue, May 26, 2020  3:06:54 PM"
#
# The underlying OOXML diff swaps the contents of ppt/theme/theme1.xml
# (the theme actually driving the slide master / visible slides -
# previously the "Integral" colour scheme) and ppt/theme/theme2.xml
# (the dormant theme only referenced by the notes master - previously
# the stock "Office Theme" colour scheme): after the edit, theme1.xml
# carries the "Office Theme" colours and theme2.xml carries the
# "Integral" colours.
#
# In other words, the presentation's live Design/Theme had its colour
# scheme swapped from the custom "Integral" palette to the default
# Office palette. We reproduce that visible effect with the
# PowerPoint object model by rewriting every entry of the active
# theme's ThemeColorScheme (the 12 DrawingML theme colours: dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) to the standard Office Theme
# RGB values.

function ToBGR($rrggbb) {
    # ThemeColorScheme.Colors(i).RGB expects the usual VBA RGB()
    # packing (0x00BBGGRR); our palette below is written as the more
    # readable 0xRRGGBB, so convert it.
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Standard "Office Theme" colours, in ThemeColorScheme.Colors index
# order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    0x000000, # dk1
    0xFFFFFF, # lt1
    0x44546A, # dk2
    0xE7E6E6, # lt2
    0x5B9BD5, # accent1
    0xED7D31, # accent2
    0xA5A5A5, # accent3
    0xFFC000, # accent4
    0x4472C4, # accent5
    0x70AD47, # accent6
    0x0563C1, # hlink
    0x954F72  # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = ToBGR $officeThemeColors[$i - 1]
}
